$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 38 (this pushes the old rows 38-85 down to 39-86,
# matching the dimension growing from A1:R85 to A1:R86).
$ws.Rows("38").Insert()

# Populate the newly inserted row 38 with the new weekly record. It mirrors
# the (now shifted-down) former row 38 except for a later date and a lower
# price (1000 instead of 1100).
$ws.Range("A38").Value = 5
$ws.Range("B38").Value = "Macroferia Regional de Talca"
$ws.Range("C38").Value = "Maule"
$ws.Range("D38").Value = 44880
$ws.Range("E38").Value = 7
$ws.Range("F38").Value = 300000000
$ws.Range("G38").Value = "Espárragos"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 3000
$ws.Range("K38").Value = 1000
$ws.Range("L38").Value = 1000
$ws.Range("M38").Value = 1000
$ws.Range("N38").Value = "$/kilo"
$ws.Range("O38").Value = "Provincia de Linares"
$ws.Range("P38").Value = 1000
$ws.Range("Q38").Value = 1
$ws.Range("R38").Value = "Hortaliza"
